# Fruta / hortaliza, semanal
# Update Fecha (D) and Volumen/Precio columns (M,N,O,P,S) for the weekly
# re-ordering of records in the "Arándano (blue)" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44455
$ws.Range("M2").Value = 160

# Row 3
$ws.Range("D3").Value = 44459
$ws.Range("M3").Value = 160

# Row 4
$ws.Range("D4").Value = 44462
$ws.Range("M4").Value = 140

# Row 5
$ws.Range("D5").Value = 44446
$ws.Range("M5").Value = 300

# Row 6
$ws.Range("D6").Value = 44463
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 13500
$ws.Range("S6").Value = 6750

# Row 7
$ws.Range("D7").Value = 44454
$ws.Range("M7").Value = 300

# Row 8
$ws.Range("D8").Value = 44445
$ws.Range("M8").Value = 160

# Row 9
$ws.Range("D9").Value = 44448
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("S9").Value = 7250

# Row 10
$ws.Range("D10").Value = 44452
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 13000
$ws.Range("P10").Value = 13500
$ws.Range("S10").Value = 6750

# Row 12
$ws.Range("D12").Value = 44466
$ws.Range("M12").Value = 160
$ws.Range("N12").Value = 13500
$ws.Range("P12").Value = 13750
$ws.Range("S12").Value = 6875
